# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values (currentAveragePrice* / Leve* profit columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H28").Value = 1284.4166
$ws.Range("I28").Value = 833.3333
$ws.Range("K28").Value = 833.3333
$ws.Range("M28").Value = -348.3333

$ws.Range("H43").Value = 849.125
$ws.Range("I43").Value = 818.2
$ws.Range("J43").Value = 900.6667
$ws.Range("K43").Value = 818.2
$ws.Range("L43").Value = 900.6667
$ws.Range("M43").Value = -749.2
$ws.Range("N43").Value = -1038.6667

$ws.Range("H70").Value = 3720.2778
$ws.Range("I70").Value = 7149.4287
$ws.Range("J70").Value = 1538.091
$ws.Range("K70").Value = 21448.2861
$ws.Range("L70").Value = 4614.272999999999
$ws.Range("M70").Value = -21178.2861
$ws.Range("N70").Value = -5154.272999999999

$ws.Range("H73").Value = 3720.2778
$ws.Range("I73").Value = 7149.4287
$ws.Range("J73").Value = 1538.091
$ws.Range("K73").Value = 21448.2861
$ws.Range("L73").Value = 4614.272999999999
$ws.Range("M73").Value = -20512.2861
$ws.Range("N73").Value = -6486.272999999999

$ws.Range("H88").Value = 791.5263
$ws.Range("I88").Value = 861.2
$ws.Range("J88").Value = 766.6429000000001
$ws.Range("K88").Value = 861.2
$ws.Range("L88").Value = 766.6429000000001
$ws.Range("M88").Value = -455.2
$ws.Range("N88").Value = -1578.6429

$ws.Range("H91").Value = 791.5263
$ws.Range("I91").Value = 861.2
$ws.Range("J91").Value = 766.6429000000001
$ws.Range("K91").Value = 861.2
$ws.Range("L91").Value = 766.6429000000001
$ws.Range("M91").Value = 542.8
$ws.Range("N91").Value = -3574.6429

$ws.Range("H99").Value = 997.4
$ws.Range("I99").Value = 476
$ws.Range("J99").Value = 1453.625
$ws.Range("K99").Value = 1428
$ws.Range("L99").Value = 4360.875
$ws.Range("M99").Value = 70
$ws.Range("N99").Value = -7356.875

$ws.Range("H125").Value = 1547.4445
$ws.Range("I125").Value = 1378.2
$ws.Range("K125").Value = 12403.8
$ws.Range("M125").Value = -9943.800000000001

$ws.Range("H127").Value = 876.8333
$ws.Range("I127").Value = 598
$ws.Range("J127").Value = 1155.6666
$ws.Range("K127").Value = 1794
$ws.Range("L127").Value = 3466.9998
$ws.Range("M127").Value = 3166
$ws.Range("N127").Value = -13386.9998

$ws.Range("H135").Value = 415.45456
$ws.Range("I135").Value = 266.85715
$ws.Range("K135").Value = 2401.71435
$ws.Range("M135").Value = 133.2856500000003

$ws.Range("H137").Value = 8098688.5
$ws.Range("I137").Value = 11365050
$ws.Range("J137").Value = 114250
$ws.Range("K137").Value = 34095150
$ws.Range("L137").Value = 342750
$ws.Range("M137").Value = -34092600
$ws.Range("N137").Value = -347850

$ws.Range("H138").Value = 3483.4866
$ws.Range("I138").Value = 753
$ws.Range("J138").Value = 5345.1816
$ws.Range("K138").Value = 2259
$ws.Range("L138").Value = 16035.5448
$ws.Range("M138").Value = 2881
$ws.Range("N138").Value = -26315.5448

$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 2082.8718
$ws.Range("I61").Value = 1600.5625
$ws.Range("J61").Value = 4287.7144
$ws.Range("K61").Value = 1600.5625
$ws.Range("L61").Value = 4287.7144
$ws.Range("M61").Value = -1388.5625
$ws.Range("N61").Value = -4711.7144

$ws.Range("H88").Value = 2828.6667
$ws.Range("I88").Value = 2753
$ws.Range("J88").Value = 2980
$ws.Range("K88").Value = 2753
$ws.Range("L88").Value = 2980
$ws.Range("M88").Value = -2347
$ws.Range("N88").Value = -3792

$ws.Range("H91").Value = 2828.6667
$ws.Range("I91").Value = 2753
$ws.Range("J91").Value = 2980
$ws.Range("K91").Value = 2753
$ws.Range("L91").Value = 2980
$ws.Range("M91").Value = -1349
$ws.Range("N91").Value = -5788

$ws.Range("H122").Value = 863293.0600000001
$ws.Range("I122").Value = 959019.3
$ws.Range("J122").Value = 1756.75
$ws.Range("K122").Value = 2877057.9
$ws.Range("L122").Value = 5270.25
$ws.Range("M122").Value = -2874607.9
$ws.Range("N122").Value = -10170.25

$ws.Range("H136").Value = 2082.8718
$ws.Range("I136").Value = 1600.5625
$ws.Range("J136").Value = 4287.7144
$ws.Range("K136").Value = 4801.6875
$ws.Range("L136").Value = 12863.1432
$ws.Range("M136").Value = -2251.6875
$ws.Range("N136").Value = -17963.1432

$ws = $wb.Worksheets("BSM")
$ws.Range("H86").Value = 2759.4075
$ws.Range("I86").Value = 2822.2354
$ws.Range("J86").Value = 2652.6
$ws.Range("K86").Value = 2822.2354
$ws.Range("L86").Value = 2652.6
$ws.Range("M86").Value = -1699.2354
$ws.Range("N86").Value = -4898.6

$ws.Range("H89").Value = 2759.4075
$ws.Range("I89").Value = 2822.2354
$ws.Range("J89").Value = 2652.6
$ws.Range("K89").Value = 14111.177
$ws.Range("L89").Value = 13263
$ws.Range("M89").Value = -8495.177
$ws.Range("N89").Value = -24495

$ws.Range("H134").Value = 3001.9412
$ws.Range("I134").Value = 2839.1538
$ws.Range("J134").Value = 3531
$ws.Range("K134").Value = 8517.4614
$ws.Range("L134").Value = 10593
$ws.Range("M134").Value = -5982.4614
$ws.Range("N134").Value = -15663

$ws = $wb.Worksheets("CRP")
$ws.Range("H16").Value = 2832.2778
$ws.Range("I16").Value = 1723.4166
$ws.Range("J16").Value = 5050
$ws.Range("K16").Value = 1723.4166
$ws.Range("L16").Value = 5050
$ws.Range("M16").Value = -1436.4166
$ws.Range("N16").Value = -5624

$ws.Range("H31").Value = 2480.8157
$ws.Range("I31").Value = 1529.5518
$ws.Range("J31").Value = 5546
$ws.Range("K31").Value = 1529.5518
$ws.Range("L31").Value = 5546
$ws.Range("M31").Value = -1234.5518
$ws.Range("N31").Value = -6136

$ws.Range("H34").Value = 2480.8157
$ws.Range("I34").Value = 1529.5518
$ws.Range("J34").Value = 5546
$ws.Range("K34").Value = 1529.5518
$ws.Range("L34").Value = 5546
$ws.Range("M34").Value = -1327.5518
$ws.Range("N34").Value = -5950

$ws.Range("H62").Value = 1002496
$ws.Range("J62").Value = 2592
$ws.Range("L62").Value = 2592
$ws.Range("N62").Value = -3840

$ws.Range("H65").Value = 1002496
$ws.Range("J65").Value = 2592
$ws.Range("L65").Value = 12960
$ws.Range("N65").Value = -19200

$ws.Range("H113").Value = 2832.2778
$ws.Range("I113").Value = 1723.4166
$ws.Range("J113").Value = 5050
$ws.Range("K113").Value = 1723.4166
$ws.Range("L113").Value = 5050
$ws.Range("M113").Value = 446.5834
$ws.Range("N113").Value = -9390

$ws.Range("H132").Value = 3524.5
$ws.Range("I132").Value = 1400
$ws.Range("J132").Value = 4232.6665
$ws.Range("K132").Value = 4200
$ws.Range("L132").Value = 12697.9995
$ws.Range("M132").Value = -1670
$ws.Range("N132").Value = -17757.9995

$ws.Range("H134").Value = 3084.5806
$ws.Range("I134").Value = 957.7143
$ws.Range("K134").Value = 2873.1429
$ws.Range("M134").Value = -338.1428999999998

$ws = $wb.Worksheets("CUL")
$ws.Range("H5").Value = 433.5
$ws.Range("I5").Value = 433.5
$ws.Range("K5").Value = 1300.5
$ws.Range("M5").Value = -1188.5

$ws.Range("H10").Value = 681
$ws.Range("I10").Value = 239.5
$ws.Range("J10").Value = 1122.5
$ws.Range("K10").Value = 718.5
$ws.Range("L10").Value = 3367.5
$ws.Range("M10").Value = -579.5
$ws.Range("N10").Value = -3645.5

$ws.Range("H131").Value = 732.3143
$ws.Range("J131").Value = 994.4286
$ws.Range("L131").Value = 2983.2858
$ws.Range("N131").Value = -13063.2858

$ws.Range("H135").Value = 433.5
$ws.Range("I135").Value = 433.5
$ws.Range("K135").Value = 3901.5
$ws.Range("M135").Value = -1366.5

$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 4157.476
$ws.Range("I80").Value = 2778.2144
$ws.Range("K80").Value = 2778.2144
$ws.Range("M80").Value = -1780.2144

$ws.Range("H83").Value = 4157.476
$ws.Range("I83").Value = 2778.2144
$ws.Range("K83").Value = 13891.072
$ws.Range("M83").Value = -8899.072

$ws.Range("H102").Value = 1722.7222
$ws.Range("I102").Value = 1726.0741
$ws.Range("J102").Value = 1712.6666
$ws.Range("K102").Value = 1726.0741
$ws.Range("L102").Value = 1712.6666
$ws.Range("M102").Value = -104.0741
$ws.Range("N102").Value = -4956.6666

$ws.Range("H113").Value = 13665.556
$ws.Range("I113").Value = 2165
$ws.Range("J113").Value = 36666.668
$ws.Range("K113").Value = 2165
$ws.Range("L113").Value = 36666.668
$ws.Range("M113").Value = 5
$ws.Range("N113").Value = -41006.668

$ws = $wb.Worksheets("LTW")
$ws.Range("H68").Value = 456086.28
$ws.Range("I68").Value = 1668133
$ws.Range("J68").Value = 1568.75
$ws.Range("K68").Value = 1668133
$ws.Range("L68").Value = 1568.75
$ws.Range("M68").Value = -1667384
$ws.Range("N68").Value = -3066.75

$ws.Range("H71").Value = 456086.28
$ws.Range("I71").Value = 1668133
$ws.Range("J71").Value = 1568.75
$ws.Range("K71").Value = 8340665
$ws.Range("L71").Value = 7843.75
$ws.Range("M71").Value = -8336921
$ws.Range("N71").Value = -15331.75

$ws.Range("H82").Value = 2193.6667
$ws.Range("I82").Value = 1790.5
$ws.Range("K82").Value = 1790.5
$ws.Range("M82").Value = -1429.5

$ws.Range("H85").Value = 2193.6667
$ws.Range("I85").Value = 1790.5
$ws.Range("K85").Value = 1790.5
$ws.Range("M85").Value = -542.5

$ws.Range("H132").Value = 2538.9312
$ws.Range("I132").Value = 2204.6
$ws.Range("J132").Value = 3281.889
$ws.Range("K132").Value = 6613.799999999999
$ws.Range("L132").Value = 9845.667000000001
$ws.Range("M132").Value = -4083.799999999999
$ws.Range("N132").Value = -14905.667

$ws = $wb.Worksheets("WVR")
$ws.Range("H62").Value = 50001250
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 50001250
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H107").Value = 970
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 1026.6666
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 3079.9998
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -6919.9998

$ws.Range("H122").Value = 2601.5925
$ws.Range("I122").Value = 1997
$ws.Range("J122").Value = 3252.6924
$ws.Range("K122").Value = 5991
$ws.Range("L122").Value = 9758.0772
$ws.Range("M122").Value = -3541
$ws.Range("N122").Value = -14658.0772

$ws.Range("H132").Value = 2049.8
$ws.Range("I132").Value = 1158.909
$ws.Range("K132").Value = 3476.727
$ws.Range("M132").Value = -946.7270000000003
